$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.75
$ws.Range("N2").Value = 7
$ws.Range("P2").Value = 3.05
$ws.Range("R2").Value = 1.83
$ws.Range("S2").Value = 2.14
$ws.Range("T2").Value = 1.53
$ws.Range("U2").Value = 2.76
$ws.Range("W2").Value = 2.32
$ws.Range("X2").Value = 30
$ws.Range("AM2").Value = 55
$ws.Range("AO2").Value = 29
$ws.Range("F3").Value = 1.9
$ws.Range("G3").Value = 1.92
$ws.Range("J3").Value = 4.1
$ws.Range("K3").Value = 4.2
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 5.1
$ws.Range("O3").Value = 1.22
$ws.Range("P3").Value = 2.4
$ws.Range("R3").Value = 1.56
$ws.Range("S3").Value = 2.66
$ws.Range("U3").Value = 2.44
$ws.Range("W3").Value = 2.08
$ws.Range("X3").Value = 26
$ws.Range("Y3").Value = 21
$ws.Range("AA3").Value = 85
$ws.Range("AB3").Value = 13
$ws.Range("AE3").Value = 46
$ws.Range("AI3").Value = 48
$ws.Range("AK3").Value = 19
$ws.Range("AM3").Value = 70
$ws.Range("AN3").Value = 9.4
$ws.Range("AO3").Value = 36
$ws.Range("L4").Value = 1.27
$ws.Range("M5").Value = 1.07
$ws.Range("O5").Value = 1.07
$ws.Range("O6").Value = 1.05
$ws.Range("G7").Value = 18.5
$ws.Range("H7").Value = 1.2
$ws.Range("I7").Value = 1.26
$ws.Range("K7").Value = 8.6
$ws.Range("M8").Value = 1.06
$ws.Range("T9").Value = 1.9
$ws.Range("I10").Value = 6.8
$ws.Range("P10").Value = 2.5
$ws.Range("R10").Value = 1.58
$ws.Range("S10").Value = 2.66
$ws.Range("T10").Value = 1.78
$ws.Range("Y10").Value = 28
$ws.Range("AK10").Value = 14
$ws.Range("AN10").Value = 6.8
$ws.Range("AO10").Value = 75
$ws.Range("T11").Value = 1.53
$ws.Range("F12").Value = 2.2
$ws.Range("G12").Value = 2.24
$ws.Range("H12").Value = 3.55
$ws.Range("I12").Value = 3.6
$ws.Range("J12").Value = 3.7
$ws.Range("K12").Value = 3.75
$ws.Range("L12").Value = 1.4
$ws.Range("P12").Value = 1.97
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 1.38
$ws.Range("V12").Value = 1.38
$ws.Range("W12").Value = 1.81
$ws.Range("Y12").Value = 14
$ws.Range("Z12").Value = 25
$ws.Range("AA12").Value = 70
$ws.Range("AB12").Value = 9.6
$ws.Range("AC12").Value = 7.8
$ws.Range("AD12").Value = 15.5
$ws.Range("AE12").Value = 44
$ws.Range("AF12").Value = 12.5
$ws.Range("AI12").Value = 55
$ws.Range("AJ12").Value = 27
$ws.Range("AK12").Value = 22
$ws.Range("AN12").Value = 17
$ws.Range("AO12").Value = 42
$ws.Range("S13").Value = 2.76
$ws.Range("U13").Value = 2
$ws.Range("AA13").Value = 11.5
$ws.Range("AB13").Value = 30
$ws.Range("AD13").Value = 10
$ws.Range("AE13").Value = 14
$ws.Range("AH13").Value = 25
